# Removed all Gear cards. Functionality of Gear cards are replaced by
# newly added skill cards. Changed skill cards cost.
#
# In this workbook, the "Gear" (展厅 / Exhibition Hall) card's effect text
# (row 7, column C) is replaced with new text describing how Gear cards
# now work via the Collector skill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "挑战开始时：从额外牌堆将《收藏家》牌洗入主牌堆。<br>`n持续：道具牌使用后横置。<br>`n挑战结束时：将所有道具牌正置，所有不在额外牌堆的《收藏家》牌放回额外牌堆。"

$ws.Range("C7").Value = $newText
$ws.Range("C7").EntireRow.RowHeight = 57

$ws.Range("C8").Select()
